# Add files via upload
#
# Adds a new "第七周周一" (week 7 / Monday) planning block to Sheet1,
# mirroring the structure of the existing weekly blocks (header row,
# column-header row, five member rows, summary row), then scrolls the
# view down so the newly added rows are in frame.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Values
# ---------------------------------------------------------------------

# Row 37: section header (merged A37:D37)
$ws.Range("A37").Value = "日期：2018.10.15 第七周周一"

# Row 38: column headers
$ws.Range("A38").Value = "组员"
$ws.Range("B38").Value = "计划内容"
$ws.Range("C38").Value = "完成情况"
$ws.Range("D38").Value = "备注"

# Row 39: 蔡旭群
$ws.Range("A39").Value = "蔡旭群"
$ws.Range("B39").Value = "完成系统管理员模块功能实体的提取"

# Rows 40-41: 区梓恒 / 游志伟 (B,C,D merged across the two rows)
$ws.Range("A40").Value = "区梓恒"
$ws.Range("B40").Value = "完成普通用户相关功能模块的实体提取"
$ws.Range("D40").Value = "由于普通用户功能模块比较多，为了保证提取实体时比较完整，此任务由两个人来完成"
$ws.Range("A41").Value = "游志伟"

# Row 42: 吴培宁
$ws.Range("A42").Value = "吴培宁"
$ws.Range("B42").Value = "完成群主角色相关模块的实体提取"

# Row 43: 郑锐芝
$ws.Range("A43").Value = "郑锐芝"
$ws.Range("B43").Value = "整合各个组员的提取的实体，绘制出完整的er图"

# Rows 44-45: summary (merged A44:D45)
$ws.Range("A44").Value = "总结："

# ---------------------------------------------------------------------
# Formatting - reuse the look of the analogous earlier weekly block
# (rows 28-36) cell by cell so the new block matches existing styles.
# ---------------------------------------------------------------------

# Header / column-header / footer rows share styles with every other
# weekly block. Paste-format one source cell into one destination cell
# at a time (rather than range-to-range) so Excel doesn't re-decompose
# the uniform 4-sided border into an "outer box only" border set.
$srcA1 = $ws.Range("A1")
foreach ($col in @("A","B","C","D")) {
    $srcA1.Copy()
    $ws.Range($col + "37").PasteSpecial(-4122)
}

$srcA2 = $ws.Range("A2")
foreach ($col in @("A","B","C","D")) {
    $srcA2.Copy()
    $ws.Range($col + "38").PasteSpecial(-4122)
}

$srcA8 = $ws.Range("A8")
foreach ($col in @("A","B","C","D")) {
    $srcA8.Copy()
    $ws.Range($col + "44").PasteSpecial(-4122)
    $srcA8.Copy()
    $ws.Range($col + "45").PasteSpecial(-4122)
}

# Row 39 (A/B bold-ish wrap text, C percent-format, D plain)
$srcA3 = $ws.Range("A3")
foreach ($col in @("A","B")) {
    $srcA3.Copy()
    $ws.Range($col + "39").PasteSpecial(-4122)
}
$ws.Range("C21").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D39").PasteSpecial(-4122)

# Rows 42-43 (plain A/B wrap text, plain C/D)
foreach ($row in @(42,43)) {
    foreach ($col in @("A","B")) {
        $srcA3.Copy()
        $ws.Range($col + $row).PasteSpecial(-4122)
    }
    foreach ($col in @("C","D")) {
        $ws.Range("C3").Copy()
        $ws.Range($col + $row).PasteSpecial(-4122)
    }
}

# Rows 40-41: A column uses the regular data style; B/C/D form a merged
# box (top cell loses its bottom border, bottom cell loses its top
# border) but keep the same 宋体/family-3 font as the rest of column A/B.
$ws.Range("A3").Copy()
$ws.Range("A40:A41").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("B40").Borders.Item(9).LineStyle = -4142

$ws.Range("A3").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("B41").Borders.Item(8).LineStyle = -4142

$ws.Range("A3").Copy()
$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("C40").Borders.Item(9).LineStyle = -4142

$ws.Range("A3").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").Borders.Item(8).LineStyle = -4142

$ws.Range("A3").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Borders.Item(9).LineStyle = -4142

$ws.Range("A3").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D41").Borders.Item(8).LineStyle = -4142

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Merges
# ---------------------------------------------------------------------
$ws.Range("A37:D37").Merge()
$ws.Range("A44:D45").Merge()
$ws.Range("B40:B41").Merge()
$ws.Range("C40:C41").Merge()
$ws.Range("D40:D41").Merge()

# ---------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(39).RowHeight = 27
$ws.Rows.Item(40).RowHeight = 29.25
$ws.Rows.Item(41).RowHeight = 24.75
$ws.Rows.Item(42).RowHeight = 19.5
$ws.Rows.Item(43).RowHeight = 27

# ---------------------------------------------------------------------
# View: scroll so the new block is in frame, move the active selection
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("H41").Select()
